$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting all existing rows (27..63) down to (28..64).
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly data point.
$ws.Cells.Item(27, 1).Value = 4
$ws.Cells.Item(27, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value = "Los Lagos"
$ws.Cells.Item(27, 4).Value = Get-Date -Year 2022 -Month 6 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(27, 5).Value = 10
$ws.Cells.Item(27, 6).Value = 100112043
$ws.Cells.Item(27, 7).Value = "Pepino dulce"
$ws.Cells.Item(27, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(27, 9).Value = "Especial"
$ws.Cells.Item(27, 10).Value = 30
$ws.Cells.Item(27, 11).Value = 21000
$ws.Cells.Item(27, 12).Value = 21000
$ws.Cells.Item(27, 13).Value = 21000
$ws.Cells.Item(27, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(27, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(27, 16).Value = 1167
$ws.Cells.Item(27, 17).Value = 18
$ws.Cells.Item(27, 18).Value = "Hortaliza"
